$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7; existing rows 7-20 shift down to 8-21.
$ws.Rows("7:7").Insert()

# Fill the newly inserted row 7 with this week's entry.
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 45274
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100103
$ws.Cells.Item(7, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(7, 9).Value = 100103003
$ws.Cells.Item(7, 10).Value = "Damasco"
$ws.Cells.Item(7, 11).Value = "Castle Brite"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 100
$ws.Cells.Item(7, 14).Value = 16000
$ws.Cells.Item(7, 15).Value = 17000
$ws.Cells.Item(7, 16).Value = 16500
$ws.Cells.Item(7, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(7, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(7, 19).Value = 917
$ws.Cells.Item(7, 20).Value = 18
